# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# 1) "ODI Batting Extra" previously stored a lot of fully-empty placeholder
#    cells (columns B-F) for rows whose scraper did not produce a value.
#    The updated scraper no longer emits those empty placeholders, so we
#    clear (and thereby drop) every currently-empty cell in the data area.
# 2) A brand new "ODI Bowling Extra" sheet is appended at the end of the
#    workbook with MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL
#    columns, mirroring the layout/style of "ODI Batting Extra".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: clean up empty placeholder cells on "ODI Batting Extra"
# ---------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$lastRow = 21
$lastCol = 6   # columns A..F

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 2; $c -le $lastCol; $c++) {
        $cell = $battingExtra.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -eq $null -or $v -eq "") {
            $cell.ClearContents()
        }
    }
}

# ---------------------------------------------------------------------
# Step 2: add the new "ODI Bowling Extra" sheet at the end of the workbook
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Copy the header formatting (bold, centered, bordered) from the sibling
# "ODI Batting Extra" sheet so both "Extra" sheets look consistent.
$battingExtra.Range("A1:C1").Copy()
$bowlingExtra.Range("A1").PasteSpecial(-4122) # xlPasteFormats

$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
$rows = @(
    @("3452", "0", $null),
    @("3506", $null, $null),
    @("3510", "0", $null),
    @("3513", "0", "10.00%"),
    @("3520", $null, $null),
    @("3580", $null, $null),
    @("3581", $null, $null),
    @("3583", "1", "20.00%"),
    @("3593", $null, $null),
    @("3596", $null, $null),
    @("3622", "0", "10.00%"),
    @("3625", "0", "10.00%"),
    @("3629", "0", "10.00%"),
    @("3655", "1", "10.00%"),
    @("3657", "1", "20.00%"),
    @("3661", $null, $null),
    @("3678", "0", "20.00%"),
    @("3680", "0", "10.00%"),
    @("3853", "0", "10.00%"),
    @("3855", "0", "10.00%")
)

$rowIndex = 2
foreach ($row in $rows) {
    $matchCode = $row[0]
    $maidenOvers = $row[1]
    $pctWickets = $row[2]

    $codeCell = $bowlingExtra.Cells.Item($rowIndex, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $matchCode

    if ($maidenOvers -ne $null) {
        $moCell = $bowlingExtra.Cells.Item($rowIndex, 2)
        $moCell.NumberFormat = "@"
        $moCell.Value = $maidenOvers
    }

    if ($pctWickets -ne $null) {
        $pctCell = $bowlingExtra.Cells.Item($rowIndex, 3)
        $pctCell.NumberFormat = "@"
        $pctCell.Value = $pctWickets
    }

    $rowIndex++
}
